# Refined metadata to be additional tab
$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# Add a new worksheet right after "data" and name it "metadata"
$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = "metadata"

# A style source cell on the "data" sheet that already carries the bold /
# bordered / centered header style we want to reuse for the metadata sheet.
$styleSource = $dataSheet.Range("B1")

# ---- Header row (row 1) ----
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

$styleSource.Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)

# ---- Data row (row 2) ----
$ws.Range("A2").Value = 0
$styleSource.Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Range("B2").Value = "Non-acute porphyrias"
$ws.Range("C2").Value = 513

# data_version "1.21" must stay textual (not be coerced into the number
# 1.21). Build it as a text formula and then convert to a plain value so
# the underlying cell keeps its text type without picking up a stray
# number-format style.
$dvCell = $ws.Range("D2")
$dvCell.Formula = '="1.21"'
$dvCell.Copy()
$dvCell.PasteSpecial(-4163)

$ws.Range("E2").Value = "2021-03-17T14:20:06.880369Z"
$ws.Range("F2").Value = "2021-10-05 14:21:50.309373"
$ws.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/513/?format=json"

$excel.CutCopyMode = $false

# ---- Refresh the "data" sheet's time_taken column (F2:F10) with the new
#      panel-query timestamps recorded for this run ----
$dataSheet.Range("F2").Value = "2021-10-05 14:21:50.313058"
$dataSheet.Range("F3").Value = "2021-10-05 14:21:50.313065"
$dataSheet.Range("F4").Value = "2021-10-05 14:21:50.313069"
$dataSheet.Range("F5").Value = "2021-10-05 14:21:50.313071"
$dataSheet.Range("F6").Value = "2021-10-05 14:21:50.313074"
$dataSheet.Range("F7").Value = "2021-10-05 14:21:50.313077"
$dataSheet.Range("F8").Value = "2021-10-05 14:21:50.313080"
$dataSheet.Range("F9").Value = "2021-10-05 14:21:50.313082"
$dataSheet.Range("F10").Value = "2021-10-05 14:21:50.313085"

$dataSheet.Activate()
